# "Use WHO tsr data" - replace the placeholder year/tx_success_pct sample
# rows on the "time_variant" sheet with the full WHO treatment-success-rate
# (tsr) time series, re-sort it, and make that sheet the active one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variant")

$years  = @(1945,1950,1995,1998,1999,2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022)
$values = @(0,50,87,83,88,91,86,94,88,89,91,86,87,89,93,90,91,89,86,87,91,90,89,92,92,83,86,86)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Header row no longer carries a row-wide custom format; drop the
# row-level formatting but keep (re-apply) the bold style on its cells.
$ws.Rows.Item(1).ClearFormats()
$ws.Range("A1:B1").Font.Bold = $true

# Match the bold formatting already used in column A.
$ws.Range("B2:B29").Font.Bold = $true

# Sort the new data by year (keeps the sortState metadata Excel records
# after running Data > Sort).
$sortRange = $ws.Range("A2:B31")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A31"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Make "time_variant" the active/selected sheet with the same selection.
$ws.Activate()
$ws.Range("E6").Select()

$wb.Save()
